$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "layer" header in D1
$ws.Range("D1").Value = "layer"

# Layer values for rows 2..44 (column D), derived from the taxon groupings
# in column B: tree layer (Picea abies) = 1, moss layer (Polytrichum
# formosum) = 9, everything else (herb layer) = 6.
$layers = @(1,6,6,6,6,6,9,
            1,6,6,6,6,6,9,
            1,6,6,6,6,9,
            1,6,6,6,6,6,6,9,
            1,6,6,6,6,6,6,9,
            1,6,6,6,6,6,9)

for ($i = 0; $i -lt $layers.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $layers[$i]
}

# Update the selection/view: select B1 (this also clears any scrolled
# topLeftCell state left over from before).
$ws.Range("B1").Select()
